$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 1 - Cumulative cases")
$ws.Range("A73").Value = 43966
